$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 23
$ws.Range("I23").Value = 6.5
$ws.Range("L23").Value = 1.36
$ws.Range("M23").Value = 3
$ws.Range("N23").Value = 2.1
$ws.Range("O23").Value = 1.7
$ws.Range("R23").Value = 2.1
$ws.Range("S23").Value = 1.67
$ws.Range("X23").Value = 15
$ws.Range("Y23").Value = 34
$ws.Range("Z23").Value = 8.5
$ws.Range("AD23").Value = 501
$ws.Range("AG23").Value = 21

# Row 24
$ws.Range("G24").Value = 1.65
$ws.Range("I24").Value = 5.5
$ws.Range("K24").Value = 10
$ws.Range("N24").Value = 2
$ws.Range("O24").Value = 1.8
$ws.Range("U24").Value = 7.5
$ws.Range("X24").Value = 15
$ws.Range("AF24").Value = 26
$ws.Range("AG24").Value = 17
$ws.Range("AJ24").Value = 41

# Row 41
$ws.Range("I41").Value = 3.9
$ws.Range("R41").Value = 2.05
$ws.Range("S41").Value = 1.7
$ws.Range("X41").Value = 21
$ws.Range("AD41").Value = 501
$ws.Range("AG41").Value = 15

# Row 57
$ws.Range("J57").Value = 1.08
$ws.Range("K57").Value = 8
$ws.Range("AD57").Value = 1000

# Row 58
$ws.Range("H58").Value = 3.5
$ws.Range("I58").Value = 5.5
$ws.Range("J58").Value = 1.13
$ws.Range("K58").Value = 6
$ws.Range("L58").Value = 1.57
$ws.Range("M58").Value = 2.25
$ws.Range("N58").Value = 2.7
$ws.Range("O58").Value = 1.44
$ws.Range("P58").Value = 1.62
$ws.Range("Q58").Value = 2.2
$ws.Range("R58").Value = 2.63
$ws.Range("S58").Value = 1.44
$ws.Range("T58").Value = 4.5
$ws.Range("U58").Value = 6
$ws.Range("V58").Value = 10
$ws.Range("X58").Value = 19
$ws.Range("Z58").Value = 6
$ws.Range("AB58").Value = 29
$ws.Range("AC58").Value = 126
$ws.Range("AE58").Value = 9.5
$ws.Range("AG58").Value = 21
$ws.Range("AJ58").Value = 67

# Row 60
$ws.Range("N60").Value = 1.8
$ws.Range("O60").Value = 2
$ws.Range("R60").Value = 2.25
$ws.Range("S60").Value = 1.57
$ws.Range("U60").Value = 6
$ws.Range("W60").Value = 8
$ws.Range("X60").Value = 12
$ws.Range("AB60").Value = 26
$ws.Range("AC60").Value = 81
$ws.Range("AH60").Value = 101

# Row 112
$ws.Range("N112").Value = 1.22
$ws.Range("O112").Value = 4.2
$ws.Range("R112").Value = 1.95
$ws.Range("S112").Value = 1.8
$ws.Range("V112").Value = 13
$ws.Range("W112").Value = 8
$ws.Range("X112").Value = 11
$ws.Range("Y112").Value = 29
$ws.Range("Z112").Value = 29
$ws.Range("AA112").Value = 19
$ws.Range("AB112").Value = 26
$ws.Range("AD112").Value = 251
$ws.Range("AE112").Value = 41
$ws.Range("AI112").Value = 101

# Row 123
$ws.Range("G123").Value = 2.3
$ws.Range("H123").Value = 3.1
$ws.Range("I123").Value = 3.2
$ws.Range("J123").Value = 1.11
$ws.Range("K123").Value = 6.5
$ws.Range("L123").Value = 1.5
$ws.Range("M123").Value = 2.5
$ws.Range("N123").Value = 2.5
$ws.Range("O123").Value = 1.5
$ws.Range("P123").Value = 1.57
$ws.Range("Q123").Value = 2.25
$ws.Range("R123").Value = 2.2
$ws.Range("S123").Value = 1.62
$ws.Range("V123").Value = 10
$ws.Range("W123").Value = 21
$ws.Range("X123").Value = 23
$ws.Range("Y123").Value = 41
$ws.Range("Z123").Value = 6.5
$ws.Range("AE123").Value = 7.5
$ws.Range("AH123").Value = 34
